# This script applies the numeric value updates to the Leve profit-tracking
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as described by the
# "chore: update Sheets via scheduled runner" commit. Each block below updates
# the currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) for one
# specific leve row on one specific sheet. Where a cell did not previously
# contain a value (or must no longer contain one), ClearContents() is used so
# the cell is emitted/omitted exactly like the target workbook.

$wb = $excel.ActiveWorkbook

# ALC sheet, row 53
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 689.0769
$ws.Range("I53").Value = 265.8
$ws.Range("J53").Value = 2100
$ws.Range("K53").Value = 265.8
$ws.Range("L53").Value = 2100
$ws.Range("M53").Value = 371.2
$ws.Range("N53").Value = -3374

# ALC sheet, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4866.6665
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 2300
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 2300
$ws.Range("M62").Value = -9376
$ws.Range("N62").Value = -3548

# ALC sheet, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 4866.6665
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 2300
$ws.Range("K65").Value = 50000
$ws.Range("L65").Value = 11500
$ws.Range("M65").Value = -46880
$ws.Range("N65").Value = -17740

# ALC sheet, row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1535.1957
$ws.Range("I129").Value = 668.53845
$ws.Range("J129").Value = 1876.6061
$ws.Range("K129").Value = 2005.61535
$ws.Range("L129").Value = 5629.8183
$ws.Range("M129").Value = 2994.38465
$ws.Range("N129").Value = -15629.8183

# ALC sheet, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1809.7428
$ws.Range("I138").Value = 1388.738
$ws.Range("J138").Value = 2441.25
$ws.Range("K138").Value = 4166.214
$ws.Range("L138").Value = 7323.75
$ws.Range("M138").Value = 973.7860000000001
$ws.Range("N138").Value = -17603.75

# ARM sheet, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3629.5
$ws.Range("I45").Value = 2262
$ws.Range("J45").Value = 4997
$ws.Range("K45").Value = 2262
$ws.Range("L45").Value = 4997
$ws.Range("M45").Value = -1885
$ws.Range("N45").Value = -5751

# BSM sheet, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2150.111
$ws.Range("I107").Value = 2168.875
$ws.Range("K107").Value = 2168.875
$ws.Range("M107").Value = -248.875

# CRP sheet, row 11
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 1300
$ws.Range("J11").Value = 1300
$ws.Range("L11").Value = 1300
$ws.Range("N11").Value = -1580

# CRP sheet, row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3129.5
$ws.Range("I16").Value = 3181.8
$ws.Range("J16").Value = 3064.125
$ws.Range("K16").Value = 3181.8
$ws.Range("L16").Value = 3064.125
$ws.Range("M16").Value = -2894.8
$ws.Range("N16").Value = -3638.125

# CRP sheet, row 43
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 17862.637
$ws.Range("J43").Value = 17862.637
$ws.Range("L43").Value = 17862.637
$ws.Range("N43").Value = -18230.637

# CRP sheet, row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 2000
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2902

# CRP sheet, row 101
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H101").Value = 17862.637
$ws.Range("J101").Value = 17862.637
$ws.Range("L101").Value = 17862.637
$ws.Range("N101").Value = -24352.637

# CRP sheet, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1835.8182
$ws.Range("I107").Value = 466.33334
$ws.Range("J107").Value = 2349.375
$ws.Range("K107").Value = 466.33334
$ws.Range("L107").Value = 2349.375
$ws.Range("M107").Value = 1453.66666
$ws.Range("N107").Value = -6189.375

# CRP sheet, row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 3129.5
$ws.Range("I113").Value = 3181.8
$ws.Range("J113").Value = 3064.125
$ws.Range("K113").Value = 3181.8
$ws.Range("L113").Value = 3064.125
$ws.Range("M113").Value = -1011.8
$ws.Range("N113").Value = -7404.125

# CRP sheet, row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 21298.777
$ws.Range("I141").Value = 5500
$ws.Range("J141").Value = 33937.8
$ws.Range("K141").Value = 5500
$ws.Range("L141").Value = 33937.8
$ws.Range("M141").Value = -320
$ws.Range("N141").Value = -44297.8

# CUL sheet, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 916.67566
$ws.Range("I131").Value = 409
$ws.Range("J131").Value = 1160.36
$ws.Range("K131").Value = 1227
$ws.Range("L131").Value = 3481.08
$ws.Range("M131").Value = 3813
$ws.Range("N131").Value = -13561.08

# CUL sheet, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1284.8235
$ws.Range("I132").Value = 781.5
$ws.Range("K132").Value = 7033.5
$ws.Range("M132").Value = -4503.5

# GSM sheet, row 3
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 4000
$ws.Range("I3").Value = 4000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 4000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -3884
$ws.Range("N3").ClearContents()

# GSM sheet, row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 5500
$ws.Range("I5").Value = 1000
$ws.Range("J5").Value = 10000
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 10000
$ws.Range("M5").Value = -888
$ws.Range("N5").Value = -10224

# GSM sheet, row 10
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 1125
$ws.Range("I10").Value = 300
$ws.Range("J10").Value = 1950
$ws.Range("K10").Value = 300
$ws.Range("L10").Value = 1950
$ws.Range("M10").Value = -131
$ws.Range("N10").Value = -2288

# GSM sheet, row 12
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 1000000
$ws.Range("I12").Value = 1000000
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1000000
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -999860
$ws.Range("N12").ClearContents()

# GSM sheet, row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 798.16
$ws.Range("I107").Value = 831.55
$ws.Range("J107").Value = 664.6
$ws.Range("K107").Value = 831.55
$ws.Range("L107").Value = 664.6
$ws.Range("M107").Value = 1088.45
$ws.Range("N107").Value = -4504.6

# LTW sheet, row 2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3754243.8
$ws.Range("I2").Value = 3000
$ws.Range("J2").Value = 5004658.5
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 5004658.5
$ws.Range("M2").Value = -2888
$ws.Range("N2").Value = -5004882.5

# WVR sheet, row 21
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 100017
$ws.Range("J21").Value = 100017
$ws.Range("L21").Value = 100017
$ws.Range("N21").Value = -100487

# WVR sheet, row 22
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 3000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 3000
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -3586

# WVR sheet, row 35
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H35").Value = 100017
$ws.Range("J35").Value = 100017
$ws.Range("L35").Value = 100017
$ws.Range("N35").Value = -100597

# WVR sheet, row 121
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 27466.666
$ws.Range("J121").Value = 27466.666
$ws.Range("L121").Value = 27466.666
$ws.Range("N121").Value = -30960.666
